# Refresh scrape snapshot to match output generated at 456a3b4.
# Values below mirror the upstream commit: ticket/interest counters bumped,
# one listing flips to sold-out ("不可售"), and three "全部类型" rows shuffle
# forward to reflect a newly scraped/resolved event.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")

# ---- 展览: refresh "想去人数"/"最低票价" counters ----
$wsExhibit.Range("F6").Value = 288
$wsExhibit.Range("F7").Value = 13289
$wsExhibit.Range("G7").Value = 238
$wsExhibit.Range("F8").Value = 78
$wsExhibit.Range("F9").Value = 128
$wsExhibit.Range("F10").Value = 316
$wsExhibit.Range("F11").Value = 5202
$wsExhibit.Range("F21").Value = 3746
$wsExhibit.Range("F24").Value = 4943
$wsExhibit.Range("F25").Value = 430
$wsExhibit.Range("F26").Value = 1978
$wsExhibit.Range("F28").Value = 287
$wsExhibit.Range("F29").Value = 7210
$wsExhibit.Range("F32").Value = 2151
$wsExhibit.Range("F33").Value = 2083
$wsExhibit.Range("F34").Value = 1313
$wsExhibit.Range("F36").Value = 1121
$wsExhibit.Range("F39").Value = 230
$wsExhibit.Range("F42").Value = 1157
$wsExhibit.Range("F45").Value = 1260
$wsExhibit.Range("F46").Value = 1894
$wsExhibit.Range("F47").Value = 86

# ---- 演出: row 6 ("缤纷国图" piano concert) flips to sold out ----
$wsShow.Range("G6").Value = "不可售"

# ---- 本地生活: refresh "想去人数" counters ----
$wsLocal.Range("F2").Value = 503
$wsLocal.Range("F3").Value = 674
$wsLocal.Range("F4").Value = 52

# ---- 全部类型: refresh "想去人数"/"最低票价" counters ----
$wsAll.Range("F7").Value = 503
$wsAll.Range("F8").Value = 674
$wsAll.Range("F9").Value = 288
$wsAll.Range("F10").Value = 13289
$wsAll.Range("G10").Value = 238
$wsAll.Range("F20").Value = 3746
$wsAll.Range("F23").Value = 4943
$wsAll.Range("F24").Value = 430
$wsAll.Range("F25").Value = 1978
$wsAll.Range("F27").Value = 288
$wsAll.Range("F28").Value = 7211
$wsAll.Range("F31").Value = 2151
$wsAll.Range("F32").Value = 2083
$wsAll.Range("F33").Value = 1313
$wsAll.Range("F35").Value = 1121
$wsAll.Range("F38").Value = 230
$wsAll.Range("F41").Value = 1157
$wsAll.Range("F45").Value = 1260
$wsAll.Range("F46").Value = 1894
$wsAll.Range("F47").Value = 86

# ---- 全部类型: rows 11-13 shuffle forward (the 07-27 piano-concert
#      listing resolved off the list; row 14 "玄色亲签...领取预约票"
#      stays put) ----
$r = 11
$wsAll.Range("B$r").NumberFormat = "@"
$wsAll.Range("B$r").Value = "2024-07-28"
$wsAll.Range("B$r").Style = "Normal"
$wsAll.Range("C$r").Value = "通州·万达动漫次元嘉年华2.0-免票活动"
$wsAll.Range("D$r").Value = "新华西街58号万达广场 北京通州万达广场"
$wsAll.Range("E$r").Value = "2024.07.28 13:00-07.28 18:00"
$wsAll.Range("F$r").Value = 316
$wsAll.Range("G$r").Value = 30
$wsAll.Range("H$r").Value = "https://show.bilibili.com/platform/detail.html?id=89056"
$wsAll.Range("I$r").Value = "//i2.hdslb.com/bfs/openplatform/202407/PVxNgbrC1720599209149.jpeg"

$r = 12
$wsAll.Range("B$r").NumberFormat = "@"
$wsAll.Range("B$r").Value = "2024-08-02"
$wsAll.Range("B$r").Style = "Normal"
$wsAll.Range("C$r").Value = "北京·IDO暑假狂欢节"
$wsAll.Range("D$r").Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
$wsAll.Range("E$r").Value = "2024.08.02 09:30-08.04 17:00"
$wsAll.Range("F$r").Value = 5202
$wsAll.Range("G$r").Value = 85
$wsAll.Range("H$r").Value = "https://show.bilibili.com/platform/detail.html?id=85556"
$wsAll.Range("I$r").Value = "//i1.hdslb.com/bfs/openplatform/202405/dzkkOQmL1716518027300.jpeg"

$r = 13
$wsAll.Range("B$r").NumberFormat = "@"
$wsAll.Range("B$r").Value = "2024-08-03"
$wsAll.Range("B$r").Style = "Normal"
$wsAll.Range("C$r").Value = "北京·玄色亲签《哑舍》周边手渡预约票"
$wsAll.Range("D$r").Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
$wsAll.Range("E$r").Value = "2024.08.03 13:00-08.03 13:30"
$wsAll.Range("F$r").Value = 14
$wsAll.Range("G$r").Value = 1
$wsAll.Range("H$r").Value = "https://show.bilibili.com/platform/detail.html?id=89722"
$wsAll.Range("I$r").Value = "//i1.hdslb.com/bfs/openplatform/202407/VzNiy0Bs1721719541445.png"

